$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1250
$ws.Range("H73").Value = 1250
$ws.Range("H88").Value = 964.7143
$ws.Range("I88").Value = 913.25
$ws.Range("J88").Value = 1033.3334
$ws.Range("K88").Value = 913.25
$ws.Range("L88").Value = 1033.3334
$ws.Range("M88").Value = -507.25
$ws.Range("N88").Value = -1845.3334
$ws.Range("H91").Value = 964.7143
$ws.Range("I91").Value = 913.25
$ws.Range("J91").Value = 1033.3334
$ws.Range("K91").Value = 913.25
$ws.Range("L91").Value = 1033.3334
$ws.Range("M91").Value = 490.75
$ws.Range("N91").Value = -3841.3334
$ws.Range("H92").Value = 83333830
$ws.Range("I92").Value = 111111560
$ws.Range("J92").Value = 666.6667
$ws.Range("K92").Value = 111111560
$ws.Range("L92").Value = 666.6667
$ws.Range("M92").Value = -111110312
$ws.Range("N92").Value = -3162.6667
$ws.Range("H111").Value = 2332.0557
$ws.Range("I111").Value = 2677
$ws.Range("J111").Value = 1124.75
$ws.Range("K111").Value = 8031
$ws.Range("L111").Value = 3374.25
$ws.Range("M111").Value = -4964
$ws.Range("N111").Value = -9508.25
$ws.Range("H116").Value = 4797.6
$ws.Range("J116").Value = 4797.6
$ws.Range("L116").Value = 4797.6
$ws.Range("N116").Value = -11681.6
$ws.Range("H133").Value = 50780
$ws.Range("J133").Value = 50780
$ws.Range("L133").Value = 50780
$ws.Range("N133").Value = -60900
$ws.Range("H137").Value = 1812.4
$ws.Range("I137").Value = 1419.0714
$ws.Range("J137").Value = 2313
$ws.Range("K137").Value = 4257.2142
$ws.Range("L137").Value = 6939
$ws.Range("M137").Value = -1707.2142
$ws.Range("N137").Value = -12039
$ws.Range("H138").Value = 2325.75
$ws.Range("J138").Value = 2347.9487
$ws.Range("L138").Value = 7043.8461
$ws.Range("N138").Value = -17323.8461
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7059.829
$ws.Range("I32").Value = 7442.788
$ws.Range("J32").Value = 5480.125
$ws.Range("K32").Value = 7442.788
$ws.Range("L32").Value = 5480.125
$ws.Range("M32").Value = -7155.788
$ws.Range("N32").Value = -6054.125
$ws.Range("H61").Value = 1555.1428
$ws.Range("I61").Value = 1378.6154
$ws.Range("K61").Value = 1378.6154
$ws.Range("M61").Value = -1166.6154
$ws.Range("H63").Value = 2234
$ws.Range("J63").Value = 3000
$ws.Range("L63").Value = 3000
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 2234
$ws.Range("J66").Value = 3000
$ws.Range("L66").Value = 15000
$ws.Range("N66").Value = -21864
$ws.Range("H102").Value = 1447.25
$ws.Range("I102").Value = 789
$ws.Range("J102").Value = 1666.6666
$ws.Range("K102").Value = 789
$ws.Range("L102").Value = 1666.6666
$ws.Range("M102").Value = 833
$ws.Range("N102").Value = -4910.6666
$ws.Range("H132").Value = 30218.566
$ws.Range("I132").Value = 1674.5853
$ws.Range("K132").Value = 5023.7559
$ws.Range("M132").Value = -2493.7559
$ws.Range("H136").Value = 1555.1428
$ws.Range("I136").Value = 1378.6154
$ws.Range("K136").Value = 4135.8462
$ws.Range("M136").Value = -1585.8462
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1461.6875
$ws.Range("I20").Value = 1577.5385
$ws.Range("K20").Value = 1577.5385
$ws.Range("M20").Value = -1330.5385
$ws.Range("H86").Value = 1308.7805
$ws.Range("I86").Value = 1168.4849
$ws.Range("K86").Value = 1168.4849
$ws.Range("M86").Value = -45.48489999999993
$ws.Range("H89").Value = 1308.7805
$ws.Range("I89").Value = 1168.4849
$ws.Range("K89").Value = 5842.424499999999
$ws.Range("M89").Value = -226.4244999999992
$ws.Range("H134").Value = 16754.334
$ws.Range("I134").Value = 30004
$ws.Range("J134").Value = 3504.6667
$ws.Range("K134").Value = 90012
$ws.Range("L134").Value = 10514.0001
$ws.Range("M134").Value = -87477
$ws.Range("N134").Value = -15584.0001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11649.238
$ws.Range("I31").Value = 17112.16
$ws.Range("J31").Value = 3615.5293
$ws.Range("K31").Value = 17112.16
$ws.Range("L31").Value = 3615.5293
$ws.Range("M31").Value = -16817.16
$ws.Range("N31").Value = -4205.5293
$ws.Range("H34").Value = 11649.238
$ws.Range("I34").Value = 17112.16
$ws.Range("J34").Value = 3615.5293
$ws.Range("K34").Value = 17112.16
$ws.Range("L34").Value = 3615.5293
$ws.Range("M34").Value = -16910.16
$ws.Range("N34").Value = -4019.5293
$ws.Range("H122").Value = 1052.5714
$ws.Range("J122").Value = 1285.3572
$ws.Range("L122").Value = 3856.0716
$ws.Range("N122").Value = -8756.071599999999
$ws.Range("H132").Value = 20737.428
$ws.Range("I132").Value = 23635.738
$ws.Range("K132").Value = 70907.21400000001
$ws.Range("M132").Value = -68377.21400000001
$ws.Range("H134").Value = 622.85364
$ws.Range("I134").Value = 509.30304
$ws.Range("J134").Value = 1091.25
$ws.Range("K134").Value = 1527.90912
$ws.Range("L134").Value = 3273.75
$ws.Range("M134").Value = 1007.09088
$ws.Range("N134").Value = -8343.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 350
$ws.Range("I26").Value = 52.5
$ws.Range("J26").Value = 573.125
$ws.Range("K26").Value = 157.5
$ws.Range("L26").Value = 1719.375
$ws.Range("M26").Value = 130.5
$ws.Range("N26").Value = -2295.375
$ws.Range("H131").Value = 170317.78
$ws.Range("J131").Value = 182669.44
$ws.Range("L131").Value = 548008.3200000001
$ws.Range("N131").Value = -558088.3200000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10212.1875
$ws.Range("J70").Value = 4698.375
$ws.Range("L70").Value = 4698.375
$ws.Range("N70").Value = -5238.375
$ws.Range("H73").Value = 10212.1875
$ws.Range("J73").Value = 4698.375
$ws.Range("L73").Value = 4698.375
$ws.Range("N73").Value = -6570.375
$ws.Range("H102").Value = 62504816
$ws.Range("I102").Value = 62504816
$ws.Range("K102").Value = 62504816
$ws.Range("M102").Value = -62503194
$ws.Range("H107").Value = 5494754
$ws.Range("I107").Value = 283.2
$ws.Range("K107").Value = 283.2
$ws.Range("M107").Value = 1636.8
$ws.Range("H126").Value = 4226.25
$ws.Range("I126").Value = 3467.3
$ws.Range("K126").Value = 10401.9
$ws.Range("M126").Value = -7931.900000000001
$ws.Range("H132").Value = 22042.852
$ws.Range("I132").Value = 3856.842
$ws.Range("J132").Value = 65234.625
$ws.Range("K132").Value = 11570.526
$ws.Range("L132").Value = 195703.875
$ws.Range("M132").Value = -9040.526
$ws.Range("N132").Value = -200763.875
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3891.889
$ws.Range("I40").Value = 3244.4
$ws.Range("J40").Value = 4701.25
$ws.Range("K40").Value = 3244.4
$ws.Range("L40").Value = 4701.25
$ws.Range("M40").Value = -3108.4
$ws.Range("N40").Value = -4973.25
$ws.Range("H132").Value = 1755.5172
$ws.Range("I132").Value = 1483.4166
$ws.Range("J132").Value = 3061.6
$ws.Range("K132").Value = 4450.2498
$ws.Range("L132").Value = 9184.799999999999
$ws.Range("M132").Value = -1920.2498
$ws.Range("N132").Value = -14244.8
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 12000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 12000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 12000
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -12946
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H81").Value = 125002216
$ws.Range("I81").Value = 2440.2
$ws.Range("J81").Value = 333335170
$ws.Range("K81").Value = 4880.4
$ws.Range("L81").Value = 666670340
$ws.Range("M81").Value = -3819.4
$ws.Range("N81").Value = -666672462
$ws.Range("H84").Value = 125002216
$ws.Range("I84").Value = 2440.2
$ws.Range("J84").Value = 333335170
$ws.Range("K84").Value = 24402
$ws.Range("L84").Value = 3333351700
$ws.Range("M84").Value = -19098
$ws.Range("N84").Value = -3333362308
$ws.Range("H100").Value = 520
$ws.Range("I100").Value = 525
$ws.Range("K100").Value = 1050
$ws.Range("M100").Value = -509
$ws.Range("H122").Value = 1147.1364
$ws.Range("I122").Value = 1020
$ws.Range("J122").Value = 1299.7
$ws.Range("K122").Value = 3060
$ws.Range("L122").Value = 3899.1
$ws.Range("M122").Value = -610
$ws.Range("N122").Value = -8799.1
$ws.Range("H132").Value = 1280.8096
$ws.Range("I132").Value = 852.8823
$ws.Range("J132").Value = 3099.5
$ws.Range("K132").Value = 2558.6469
$ws.Range("L132").Value = 9298.5
$ws.Range("M132").Value = -28.64689999999973
$ws.Range("N132").Value = -14358.5
